# Generate Report for Handback
# Fills in the "42de5674-7117-40fc-803a-cacb4a86b2c0" (e2e) row's handback
# columns on the per-language status sheets (zh-cn, de-de), now that its
# handback has been processed: latest target/handback xlf file, handback
# datetime, a "stale handback" error message in the Error Detail column,
# and a hyperlink on the Latest Target File cell pointing at the source .md.

$wb = $excel.ActiveWorkbook

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/831881d80340683058c235074034a851efee39f3/e2e/42de5674-7117-40fc-803a-cacb4a86b2c0.md"
$mdDisplay = "42de5674-7117-40fc-803a-cacb4a86b2c0.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0dad8075df5deaf7fd83db3cf90781cf0a8b9dd7/e2e/42de5674-7117-40fc-803a-cacb4a86b2c0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/831881d80340683058c235074034a851efee39f3/e2e/42de5674-7117-40fc-803a-cacb4a86b2c0.md."

# --- zh-cn sheet, row 7 ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value = "42de5674-7117-40fc-803a-cacb4a86b2c0.1016db0d1fa901f4f45dda943133feabbda6e9d5.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-22 05:04:28"
$wsZh.Range("P7").Value = $errorDetail

$iZh = $wsZh.Range("I7")
$iZh.Value = $mdDisplay
$wsZh.Hyperlinks.Add($iZh, $mdUrl, "", "", $mdDisplay)

# --- de-de sheet, row 7 ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value = "42de5674-7117-40fc-803a-cacb4a86b2c0.1016db0d1fa901f4f45dda943133feabbda6e9d5.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-22 05:04:34"
$wsDe.Range("P7").Value = $errorDetail

$iDe = $wsDe.Range("I7")
$iDe.Value = $mdDisplay
$wsDe.Hyperlinks.Add($iDe, $mdUrl, "", "", $mdDisplay)
